$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New regenerated sval data (rows 2-8), columns B:G
# Column F ("Win") values are unchanged.

$data = @{
    2 = @{ B = 0.6545652718822623;  C = 1.626987699542094;  D = 0.7210945179870265;  E = 0.5333859586016987;  G = 3.536033448013082  }
    3 = @{ B = 0.6545652718822623;  C = 1.626987699542094;  D = 0.1496068669990043;  E = 13.86384647080068;   G = 16.29500630922404  }
    4 = @{ B = 0.04172184405617529; C = 0.3048912486333797; D = 0.7210945179870265;  E = 0.5333859586016987;  G = 1.60109356927828   }
    5 = @{ B = 3.272327238179451;   C = 1.626987699542094;  D = 0.1496068669990043;  E = 0.5333859586016987;  G = 5.582307763322248  }
    6 = @{ B = 3.272327238179451;   C = 1.626987699542094;  D = 18.71679738969934;   E = 0.5333859586016987;  G = 24.14949828602258  }
    7 = @{ B = 0.6545652718822623;  C = 1.626987699542094;  D = 0.7210945179870265;  E = 0.5333859586016987;  G = 3.536033448013082  }
    8 = @{ B = 3.272327238179451;   C = 1.626987699542094;  D = 18.71679738969934;   E = 13.86384647080068;   G = 37.47995879822157  }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("B$row").Value = $vals.B
    $ws.Range("C$row").Value = $vals.C
    $ws.Range("D$row").Value = $vals.D
    $ws.Range("E$row").Value = $vals.E
    $ws.Range("G$row").Value = $vals.G
}
